$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Break the two old merges so every cell can be written independently.
# ---------------------------------------------------------------------------
$ws.Range("B3:B4").UnMerge()
$ws.Range("B5:B6").UnMerge()

# ---------------------------------------------------------------------------
# 2. Rewrite the B/C/D table for rows 3-17 with the new content.
#    (columns C/D keep the plain "centered" header style untouched; only
#    column B carries per-row alignment variations)
# ---------------------------------------------------------------------------

# Row 3 - News / GET / news/trending  (top of B3:B4 merge)
$ws.Range("B3").Value = "News"
$ws.Range("C3").Value = "GET"
$ws.Range("D3").Value = "news/trending"

# Row 4 - (blank) / GET / news/id  (bottom of B3:B4 merge)
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "GET"
$ws.Range("D4").Value = "news/id"

# Row 5 - Stats / GET / stats  (standalone, vertical-center only)
$ws.Range("B5").Value = "Stats"
$ws.Range("C5").Value = "GET"
$ws.Range("D5").Value = "stats"

# Row 6 - Market / GET / markets  (top of new B6:B9 merge)
$ws.Range("B6").Value = "Market"
$ws.Range("C6").Value = "GET"
$ws.Range("D6").Value = "markets"

# Row 7 - (blank) / GET / markets/id
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "GET"
$ws.Range("D7").Value = "markets/id"

# Row 8 - (blank) / GET / markets/id/categories
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "GET"
$ws.Range("D8").Value = "markets/id/categories"

# Row 9 - (blank) / GET / markets/id/categories/categoryID
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "GET"
$ws.Range("D9").Value = "markets/id/categories/categoryID"

# Row 10 - Product In Market / GET / productinmarket/marketID/productID (wrap)
$ws.Range("B10").Value = "Product In Market"
$ws.Range("C10").Value = "GET"
$ws.Range("D10").Value = "productinmarket/marketID/productID"

# Row 11 - Product / GET / product/meatID (standalone, vertical-center only)
$ws.Range("B11").Value = "Product"
$ws.Range("C11").Value = "GET"
$ws.Range("D11").Value = "product/meatID"

# Row 12 - Feedback / GET / feedback/meatID (top of new B12:B13 merge)
$ws.Range("B12").Value = "Feedback"
$ws.Range("C12").Value = "GET"
$ws.Range("D12").Value = "feedback/meatID"

# Row 13 - (blank) / POST / feedback/meatID
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "POST"
$ws.Range("D13").Value = "feedback/meatID"

# Row 14 - Saved / POST / saved/userID/productID (top of new B14:B15 merge)
$ws.Range("B14").Value = "Saved"
$ws.Range("C14").Value = "POST"
$ws.Range("D14").Value = "saved/userID/productID"

# Row 15 - (blank) / GET / saved/userID
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "GET"
$ws.Range("D15").Value = "saved/userID"

# Row 16 - User (standalone label, no C/D)
$ws.Range("B16").Value = "User"

# Row 17 - App Config (standalone label, no C/D)
$ws.Range("B17").Value = "App Config"

# ---------------------------------------------------------------------------
# 3. Re-create the merges for the new layout.
# ---------------------------------------------------------------------------
$ws.Range("B3:B4").Merge()
$ws.Range("B6:B9").Merge()
$ws.Range("B12:B13").Merge()
$ws.Range("B14:B15").Merge()

# ---------------------------------------------------------------------------
# 4. Apply alignment per the new style set:
#    - "center"  -> horizontal + vertical centered (rows 3,4,6,7,8,9,12,13,14,15)
#    - "vcenter" -> vertical centered only (rows 5, 11)
#    - "wrap"    -> horizontal + vertical centered + wrap text (row 10)
#    - rows 16/17 reuse the original centered header-label style
# ---------------------------------------------------------------------------
$centerRows = @(3,4,6,7,8,9,12,13,14,15)
foreach ($r in $centerRows) {
    $cell = $ws.Range("B$r")
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

$vcenterRows = @(5,11)
foreach ($r in $vcenterRows) {
    $cell = $ws.Range("B$r")
    $cell.HorizontalAlignment = -4142
    $cell.VerticalAlignment = -4108
}

$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("B10").VerticalAlignment = -4108
$ws.Range("B10").WrapText = $true

$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("B17").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Row height for the wrapped row and new column width for D.
# ---------------------------------------------------------------------------
$ws.Rows(10).RowHeight = 28.8
$ws.Columns("D").ColumnWidth = 39.6

# ---------------------------------------------------------------------------
# 6. Move the active selection to the new last cell, like the source file.
# ---------------------------------------------------------------------------
$ws.Range("D17").Select()
